# Apply the cryptos-list refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-sensed as a number by Excel
# (e.g. "0.9993", "7.780") are pinned to Text format first, written, then
# restored to the Normal style so no stray number-format survives on them -
# matching how the feed always stores these as plain text.

$ws.Range("D2").Value = "29.484.78"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.877.99"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7175"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07890"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3082"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08205"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "1.863.68"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.273"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7256"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "29.464.93"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.848"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007859"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("D21").Value = "2.117.96"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.780"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.94%  "
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1481"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.68%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.992"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.950"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.362"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.481"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.350"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.104"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05249"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7205"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.08%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.672"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01858"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.707"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").Value = "1.182.68"
$ws.Range("E41").Value = "  +4.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9125"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.986"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4318"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9995"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5337"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.778"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.897"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.220"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.71%  "
